# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
# D = Price, E = Volume(1h); cells are plain text in the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.671.87"
$ws.Range("E2").Value = "  -1.15%  "

$ws.Range("D3").Value = "3.784.82"
$ws.Range("E3").Value = "  -1.85%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.63"
$ws.Range("E5").Value = "  -0.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.39"
$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("D7").Value = "3.783.95"
$ws.Range("E7").Value = "  -1.84%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  -0.89%  "

$ws.Range("E10").Value = "  -0.32%  "

$ws.Range("E11").Value = "  +0.38%  "

$ws.Range("E12").Value = "  -0.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000277"
$ws.Range("E13").Value = "  +4.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.51"
$ws.Range("E14").Value = "  -1.27%  "

$ws.Range("D15").Value = "4.418.14"
$ws.Range("E15").Value = "  -1.87%  "

$ws.Range("D16").Value = "3.783.26"
$ws.Range("E16").Value = "  -1.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.58"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").Value = "67.710.97"
$ws.Range("E18").Value = "  -1.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.19"
$ws.Range("E19").Value = "  -2.33%  "

$ws.Range("E20").Value = "  +0.91%  "

$ws.Range("E21").Value = "  -7.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.87"
$ws.Range("E22").Value = "  -0.37%  "

$ws.Range("E23").Value = "  -1.70%  "

$ws.Range("E24").Value = "  -7.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.86"
$ws.Range("E25").Value = "  +0.48%  "

$ws.Range("E26").Value = "  -1.05%  "

$ws.Range("E27").Value = "  +0.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.31"
$ws.Range("E28").Value = "  +0.97%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.91"
$ws.Range("E30").Value = "  -1.50%  "

$ws.Range("D31").Value = "3.934.19"
$ws.Range("E31").Value = "  -1.84%  "

$ws.Range("E32").Value = "  -0.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.56"
$ws.Range("E33").Value = "  -2.84%  "

$ws.Range("E34").Value = "  -3.40%  "

$ws.Range("E35").Value = "  -1.98%  "

$ws.Range("D36").Value = "3.747.65"
$ws.Range("E36").Value = "  -1.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.75"
$ws.Range("E37").Value = "  +0.69%  "

$ws.Range("E38").Value = "  -0.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -1.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.137"
$ws.Range("E40").Value = "  -2.21%  "

$ws.Range("E41").Value = "  -2.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.16%  "

$ws.Range("E43").Value = "  -0.63%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.67"
$ws.Range("E45").Value = "  -0.28%  "

$ws.Range("E46").Value = "  -2.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.88"
$ws.Range("E47").Value = "  -2.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "395.05"
$ws.Range("E48").Value = "  -5.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000269"
$ws.Range("E49").Value = "  -7.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.84"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.23"
$ws.Range("E51").Value = "  +3.41%  "
